# Weekly fruit/veg price update: a new daily-price record was collected for
# "Feria Lagunitas de Puerto Montt - Acelga" and inserted as the new row 302
# (sorted position), pushing the former rows 302-339 down to 303-340.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 302; Excel shifts rows 302:339 down to 303:340,
# carrying their values/formatting with them (matches the diff exactly).
$ws.Rows.Item(302).Insert()

# Populate the newly inserted row with the new observation.
$ws.Cells.Item(302, 1).Value  = 4
$ws.Cells.Item(302, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(302, 3).Value  = "Los Lagos"
$ws.Cells.Item(302, 4).Value  = 45131
$ws.Cells.Item(302, 5).Value  = 10
$ws.Cells.Item(302, 6).Value  = 100112009
$ws.Cells.Item(302, 7).Value  = "Acelga"
$ws.Cells.Item(302, 8).Value  = "Sin especificar"
$ws.Cells.Item(302, 9).Value  = "Primera"
$ws.Cells.Item(302, 10).Value = 40
$ws.Cells.Item(302, 11).Value = 10000
$ws.Cells.Item(302, 12).Value = 10000
$ws.Cells.Item(302, 13).Value = 10000
$ws.Cells.Item(302, 14).Value = "$/docena de atados (12 kilos)"
$ws.Cells.Item(302, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(302, 16).Value = 833
$ws.Cells.Item(302, 17).Value = 12
$ws.Cells.Item(302, 18).Value = "Hortaliza"
